# Aggiornato file base.xlsx per rispettare nuovi parametri accettati
#
# - DatiGenerali!A1 becomes the single header cell "SERIAL NUMBER:" (the old
#   two-row NOME/MODELLO header is collapsed to one row).
# - DatiGenerali!A2 ("MODELLO") is removed entirely.
# - Column A on DatiGenerali is widened to fit the new, longer label.
# - The DatiGenerali sheet becomes the active/selected tab (with B3 selected),
#   taking over from the Grafico sheet which was previously active.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatiGenerali")

# Replace the header text and drop the old second row ("MODELLO").
$ws.Range("A1").Value = "SERIAL NUMBER:"
$ws.Rows.Item(2).Delete() | Out-Null

# Widen column A to fit the new, longer label text (mirrors Excel's
# "AutoFit" best-fit sizing for the "SERIAL NUMBER:" header).
$ws.Columns.Item(1).ColumnWidth = 14.86

# Make DatiGenerali the active sheet/tab with B3 selected.
$ws.Activate()
$ws.Range("B3").Select() | Out-Null
